$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to text so values like "1.001" or
# "0.000007712" are not auto-converted to numbers, then restore the
# default (unstyled) cell style so the saved XML matches the original
# (no explicit "s" attribute) while keeping the text cell type.
$ws.Range("D2:E51").NumberFormat = "@"

# --- Price / Volume(1h) updates ---
$ws.Range("D2").Value = "30.233.32"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").Value = "1.876.03"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "235.61"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.4827"
$ws.Range("E7").Value = "  -1.68%  "
$ws.Range("D8").Value = "0.2871"
$ws.Range("E8").Value = "  -2.91%  "
$ws.Range("D9").Value = "0.06566"
$ws.Range("E9").Value = "  -2.87%  "
$ws.Range("D10").Value = "1.873.59"
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("D11").Value = "16.67"
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("D12").Value = "0.07312"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "5.134"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "86.75"
$ws.Range("E14").Value = "  -3.41%  "
$ws.Range("D15").Value = "0.6517"
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("D16").Value = "30.202.94"
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("D17").Value = "13.30"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "0.000007712"
$ws.Range("E19").Value = "  -2.78%  "
$ws.Range("D20").Value = "2.133.88"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").Value = "5.308"
$ws.Range("E21").Value = "  +3.21%  "
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "192.13"
$ws.Range("E23").Value = "  -7.52%  "
$ws.Range("D24").Value = "6.094"
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("D25").Value = "9.246"
$ws.Range("E25").Value = "  -4.27%  "
$ws.Range("D26").Value = "161.65"
$ws.Range("E26").Value = "  +2.79%  "
$ws.Range("D27").Value = "17.94"
$ws.Range("E27").Value = "  -4.90%  "
$ws.Range("D28").Value = "1.902"
$ws.Range("E28").Value = "  -3.34%  "
$ws.Range("D29").Value = "1.436"
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("D30").Value = "4.248"
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("D31").Value = "0.09075"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").Value = "3.997"
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("D33").Value = "0.05059"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("D34").Value = "0.7113"
$ws.Range("E34").Value = "  -5.25%  "
$ws.Range("D35").Value = "1.091"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").Value = "2.701"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("D37").Value = "0.01782"
$ws.Range("E37").Value = "  -3.19%  "
$ws.Range("D38").Value = "2.632"
$ws.Range("E38").Value = "  -3.64%  "
$ws.Range("D39").Value = "0.9187"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("D40").Value = "2.034"
$ws.Range("E40").Value = "  -2.90%  "
$ws.Range("D41").Value = "105.62"
$ws.Range("E41").Value = "  -1.14%  "
$ws.Range("D45").Value = "7.354"
$ws.Range("E45").Value = "  -4.58%  "
$ws.Range("D46").Value = "0.1310"
$ws.Range("E46").Value = "  -5.31%  "
$ws.Range("D47").Value = "64.80"
$ws.Range("E47").Value = "  -2.16%  "
$ws.Range("D48").Value = "8.887"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("D51").Value = "0.3797"
$ws.Range("E51").Value = "  -7.08%  "

# --- Rows that were reordered (Coin/Link/Price/Volume all changed) ---
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.774"
$ws.Range("E42").Value = "  -1.46%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.4239"
$ws.Range("E44").Value = "  -5.66%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.05745"
$ws.Range("E49").Value = "  -3.00%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "33.64"
$ws.Range("E50").Value = "  -4.27%  "

# Restore default styling on the Price/Volume columns (drop the temporary
# text NumberFormat) so cells keep their original unstyled appearance.
$ws.Range("D2:E51").Style = "Normal"
